$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu results table: B column setpoint changed from 1.05 to 1.02 pu
# (per-unit bus voltage magnitudes), with cascading recomputation of the
# power-flow results in columns C:F and I:N for rows 2-25 (case with 380 kV done).

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.074594173660594
$ws.Cells.Item(2, 4).Value = 1.083821771554869
$ws.Cells.Item(2, 5).Value = 1.068266598740452
$ws.Cells.Item(2, 6).Value = 1.090032741912005
$ws.Cells.Item(2, 9).Value = 1.052528360026868
$ws.Cells.Item(2, 10).Value = 1.079502870546629
$ws.Cells.Item(2, 11).Value = 1.086486391214367
$ws.Cells.Item(2, 12).Value = 1.070972240506091
$ws.Cells.Item(2, 13).Value = 1.092681334920577
$ws.Cells.Item(2, 14).Value = 1.081035888443515
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.0764186059026
$ws.Cells.Item(3, 4).Value = 1.085585021524453
$ws.Cells.Item(3, 5).Value = 1.069870332669996
$ws.Cells.Item(3, 6).Value = 1.091836470220521
$ws.Cells.Item(3, 9).Value = 1.053050449901971
$ws.Cells.Item(3, 10).Value = 1.080982358630088
$ws.Cells.Item(3, 11).Value = 1.088066257018206
$ws.Cells.Item(3, 12).Value = 1.072389962631716
$ws.Cells.Item(3, 13).Value = 1.0943027620822
$ws.Cells.Item(3, 14).Value = 1.082517477569752
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.077595459414612
$ws.Cells.Item(4, 4).Value = 1.08672261200628
$ws.Cells.Item(4, 5).Value = 1.070904352765943
$ws.Cells.Item(4, 6).Value = 1.09300037281643
$ws.Cells.Item(4, 9).Value = 1.05338507294997
$ws.Cells.Item(4, 10).Value = 1.081935700019168
$ws.Cells.Item(4, 11).Value = 1.089084727623327
$ws.Cells.Item(4, 12).Value = 1.073303138594765
$ws.Cells.Item(4, 13).Value = 1.095348250238596
$ws.Cells.Item(4, 14).Value = 1.083472172812955
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.078089347320117
$ws.Cells.Item(5, 4).Value = 1.087200070359879
$ws.Cells.Item(5, 5).Value = 1.071338184822814
$ws.Cells.Item(5, 6).Value = 1.093488921633529
$ws.Cells.Item(5, 9).Value = 1.053524987826899
$ws.Cells.Item(5, 10).Value = 1.082335546001151
$ws.Cells.Item(5, 11).Value = 1.089511997214589
$ws.Cells.Item(5, 12).Value = 1.073686051365006
$ws.Cells.Item(5, 13).Value = 1.095786907961657
$ws.Cells.Item(5, 14).Value = 1.083872586622095
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.078172223306673
$ws.Cells.Item(6, 4).Value = 1.087280192204829
$ws.Cells.Item(6, 5).Value = 1.071410976594119
$ws.Cells.Item(6, 6).Value = 1.093570907311991
$ws.Cells.Item(6, 9).Value = 1.053548435729754
$ws.Cells.Item(6, 10).Value = 1.082402627303049
$ws.Cells.Item(6, 11).Value = 1.089583685618147
$ws.Cells.Item(6, 12).Value = 1.073750286708385
$ws.Cells.Item(6, 13).Value = 1.095860510242871
$ws.Cells.Item(6, 14).Value = 1.083939763187136
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.077602062127818
$ws.Cells.Item(7, 4).Value = 1.08672899488782
$ws.Cells.Item(7, 5).Value = 1.070910153047109
$ws.Cells.Item(7, 6).Value = 1.093006903775726
$ws.Cells.Item(7, 9).Value = 1.053386945479504
$ws.Cells.Item(7, 10).Value = 1.081941046448107
$ws.Cells.Item(7, 11).Value = 1.089090440316909
$ws.Cells.Item(7, 12).Value = 1.07330825894873
$ws.Cells.Item(7, 13).Value = 1.095354114983074
$ws.Cells.Item(7, 14).Value = 1.083477526834435
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.075211521598723
$ws.Cells.Item(8, 4).Value = 1.084418374302491
$ws.Cells.Item(8, 5).Value = 1.068809364847135
$ws.Cells.Item(8, 6).Value = 1.090642999758643
$ws.Cells.Item(8, 9).Value = 1.052705470027116
$ws.Cells.Item(8, 10).Value = 1.080003705302913
$ws.Cells.Item(8, 11).Value = 1.087021112600376
$ws.Cells.Item(8, 12).Value = 1.071452242256818
$ws.Cells.Item(8, 13).Value = 1.093230076447098
$ws.Cells.Item(8, 14).Value = 1.0815374344426
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.070970081968769
$ws.Cells.Item(9, 4).Value = 1.080320324340813
$ws.Cells.Item(9, 5).Value = 1.065078401930835
$ws.Cells.Item(9, 6).Value = 1.086451969643309
$ws.Cells.Item(9, 9).Value = 1.051479795089455
$ws.Cells.Item(9, 10).Value = 1.076558627918822
$ws.Cells.Item(9, 11).Value = 1.083344793120999
$ws.Cells.Item(9, 12).Value = 1.068148951442094
$ws.Cells.Item(9, 13).Value = 1.089458297690197
$ws.Cells.Item(9, 14).Value = 1.078087464653468
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.068121765801143
$ws.Cells.Item(10, 4).Value = 1.077569412177621
$ws.Cells.Item(10, 5).Value = 1.062570476707119
$ws.Cells.Item(10, 6).Value = 1.083639674869307
$ws.Cells.Item(10, 9).Value = 1.050645577641432
$ws.Cells.Item(10, 10).Value = 1.074239918608657
$ws.Cells.Item(10, 11).Value = 1.080872779237148
$ws.Cells.Item(10, 12).Value = 1.065923752983362
$ws.Cells.Item(10, 13).Value = 1.086923273397528
$ws.Cells.Item(10, 14).Value = 1.075765462510122
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.066883239387068
$ws.Cells.Item(11, 4).Value = 1.076373512200584
$ws.Cells.Item(11, 5).Value = 1.061479392702046
$ws.Cells.Item(11, 6).Value = 1.082417336794687
$ws.Cells.Item(11, 9).Value = 1.050280204451295
$ws.Cells.Item(11, 10).Value = 1.073230453794646
$ws.Cells.Item(11, 11).Value = 1.079797125944207
$ws.Cells.Item(11, 12).Value = 1.064954540659411
$ws.Cells.Item(11, 13).Value = 1.085820483214738
$ws.Cells.Item(11, 14).Value = 1.074754564140288
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.06642239317299
$ws.Cells.Item(12, 4).Value = 1.07592856864475
$ws.Cells.Item(12, 5).Value = 1.061073323574622
$ws.Cells.Item(12, 6).Value = 1.081962593955228
$ws.Cells.Item(12, 9).Value = 1.050143856857571
$ws.Cells.Item(12, 10).Value = 1.07285465578814
$ws.Cells.Item(12, 11).Value = 1.079396770860588
$ws.Cells.Item(12, 12).Value = 1.064593658565384
$ws.Cells.Item(12, 13).Value = 1.085410070206688
$ws.Cells.Item(12, 14).Value = 1.074378232457508
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.066521282932977
$ws.Cells.Item(13, 4).Value = 1.076024044066763
$ws.Cells.Item(13, 5).Value = 1.061160462947627
$ws.Cells.Item(13, 6).Value = 1.082060170409294
$ws.Cells.Item(13, 9).Value = 1.050173132595709
$ws.Cells.Item(13, 10).Value = 1.07293530396178
$ws.Cells.Item(13, 11).Value = 1.079482685344054
$ws.Cells.Item(13, 12).Value = 1.064671108836102
$ws.Cells.Item(13, 13).Value = 1.085498141159633
$ws.Cells.Item(13, 14).Value = 1.074458995160806
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.06684516219791
$ws.Cells.Item(14, 4).Value = 1.076336748106389
$ws.Cells.Item(14, 5).Value = 1.061445843175429
$ws.Cells.Item(14, 6).Value = 1.082379762261232
$ws.Cells.Item(14, 9).Value = 1.050268946850311
$ws.Cells.Item(14, 10).Value = 1.073199407418028
$ws.Cells.Item(14, 11).Value = 1.079764049099805
$ws.Cells.Item(14, 12).Value = 1.064924727963212
$ws.Cells.Item(14, 13).Value = 1.085786574532874
$ws.Cells.Item(14, 14).Value = 1.074723473674255
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.067044607885473
$ws.Cells.Item(15, 4).Value = 1.076529317603614
$ws.Cells.Item(15, 5).Value = 1.061621569798596
$ws.Cells.Item(15, 6).Value = 1.08257657834469
$ws.Cells.Item(15, 9).Value = 1.050327897224128
$ws.Cells.Item(15, 10).Value = 1.073362018672701
$ws.Cells.Item(15, 11).Value = 1.079937298717124
$ws.Cells.Item(15, 12).Value = 1.065080874756758
$ws.Cells.Item(15, 13).Value = 1.085964182875097
$ws.Cells.Item(15, 14).Value = 1.07488631585556
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.06820385147614
$ws.Cells.Item(16, 4).Value = 1.077648678500345
$ws.Cells.Item(16, 5).Value = 1.062642778422314
$ws.Cells.Item(16, 6).Value = 1.083720698756668
$ws.Cells.Item(16, 9).Value = 1.050669738108289
$ws.Cells.Item(16, 10).Value = 1.074306797033584
$ws.Cells.Item(16, 11).Value = 1.080944054381472
$ws.Cells.Item(16, 12).Value = 1.065987954974569
$ws.Cells.Item(16, 13).Value = 1.086996352611633
$ws.Cells.Item(16, 14).Value = 1.075832435910084
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.068929609499791
$ws.Cells.Item(17, 4).Value = 1.078349540808471
$ws.Cells.Item(17, 5).Value = 1.063281966321716
$ws.Cells.Item(17, 6).Value = 1.084437129674088
$ws.Cells.Item(17, 9).Value = 1.050883048766964
$ws.Cells.Item(17, 10).Value = 1.074897959327541
$ws.Cells.Item(17, 11).Value = 1.081574144630133
$ws.Cells.Item(17, 12).Value = 1.06655540651707
$ws.Cells.Item(17, 13).Value = 1.087642423685485
$ws.Cells.Item(17, 14).Value = 1.076424437722305
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.069352433009192
$ws.Cells.Item(18, 4).Value = 1.078757886356447
$ws.Cells.Item(18, 5).Value = 1.063654299663103
$ws.Cells.Item(18, 6).Value = 1.084854569742944
$ws.Cells.Item(18, 9).Value = 1.051007069249086
$ws.Cells.Item(18, 10).Value = 1.075242249916445
$ws.Cells.Item(18, 11).Value = 1.081941159951641
$ws.Cells.Item(18, 12).Value = 1.066885844020135
$ws.Cells.Item(18, 13).Value = 1.088018774549843
$ws.Cells.Item(18, 14).Value = 1.076769217243336
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.069496520977566
$ws.Cells.Item(19, 4).Value = 1.078897045004925
$ws.Cells.Item(19, 5).Value = 1.063781172506561
$ws.Cells.Item(19, 6).Value = 1.084996831702588
$ws.Cells.Item(19, 9).Value = 1.051049289419228
$ws.Cells.Item(19, 10).Value = 1.075359555746102
$ws.Cells.Item(19, 11).Value = 1.082066217525953
$ws.Cells.Item(19, 12).Value = 1.066998422398927
$ws.Cells.Item(19, 13).Value = 1.088147017716173
$ws.Cells.Item(19, 14).Value = 1.076886689660727
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.068851794299097
$ws.Cells.Item(20, 4).Value = 1.078274392192586
$ws.Cells.Item(20, 5).Value = 1.063213438767849
$ws.Cells.Item(20, 6).Value = 1.084360309286837
$ws.Cells.Item(20, 9).Value = 1.05086020398526
$ws.Cells.Item(20, 10).Value = 1.074834587585619
$ws.Cells.Item(20, 11).Value = 1.081506594356954
$ws.Cells.Item(20, 12).Value = 1.066494581092889
$ws.Cells.Item(20, 13).Value = 1.087573157353642
$ws.Cells.Item(20, 14).Value = 1.076360975985241
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.066749810163001
$ws.Cells.Item(21, 4).Value = 1.076244684969857
$ws.Cells.Item(21, 5).Value = 1.061361827837585
$ws.Cells.Item(21, 6).Value = 1.082285670257304
$ws.Cells.Item(21, 9).Value = 1.050240749443897
$ws.Cells.Item(21, 10).Value = 1.073121658816947
$ws.Cells.Item(21, 11).Value = 1.079681216996675
$ws.Cells.Item(21, 12).Value = 1.064850067718735
$ws.Cells.Item(21, 13).Value = 1.085701659999712
$ws.Cells.Item(21, 14).Value = 1.074645614661241
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.065423555441909
$ws.Cells.Item(22, 4).Value = 1.074964275435175
$ws.Cells.Item(22, 5).Value = 1.060193054322572
$ws.Cells.Item(22, 6).Value = 1.080977131739033
$ws.Cells.Item(22, 9).Value = 1.049847614789612
$ws.Cells.Item(22, 10).Value = 1.072039815368257
$ws.Cells.Item(22, 11).Value = 1.078528835825395
$ws.Cells.Item(22, 12).Value = 1.063811033538227
$ws.Cells.Item(22, 13).Value = 1.084520408123553
$ws.Cells.Item(22, 14).Value = 1.073562234870766
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.06612707709099
$ws.Cells.Item(23, 4).Value = 1.075643455017329
$ws.Cells.Item(23, 5).Value = 1.060813085428577
$ws.Cells.Item(23, 6).Value = 1.081671211577751
$ws.Cells.Item(23, 9).Value = 1.050056372439718
$ws.Cells.Item(23, 10).Value = 1.072613787950666
$ws.Cells.Item(23, 11).Value = 1.079140186587963
$ws.Cells.Item(23, 12).Value = 1.064362331611017
$ws.Cells.Item(23, 13).Value = 1.085147051833085
$ws.Cells.Item(23, 14).Value = 1.074137022560077
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.068886957179692
$ws.Cells.Item(24, 4).Value = 1.078308350024611
$ws.Cells.Item(24, 5).Value = 1.063244404946123
$ws.Cells.Item(24, 6).Value = 1.084395022478318
$ws.Cells.Item(24, 9).Value = 1.050870527794626
$ws.Cells.Item(24, 10).Value = 1.074863224167834
$ws.Cells.Item(24, 11).Value = 1.081537118979742
$ws.Cells.Item(24, 12).Value = 1.066522067174793
$ws.Cells.Item(24, 13).Value = 1.087604457345609
$ws.Cells.Item(24, 14).Value = 1.076389653234687
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.0720701543169
$ws.Cells.Item(25, 4).Value = 1.081383014450154
$ws.Cells.Item(25, 5).Value = 1.066046501521448
$ws.Cells.Item(25, 6).Value = 1.087538590190914
$ws.Cells.Item(25, 9).Value = 1.051799645340678
$ws.Cells.Item(25, 10).Value = 1.077453065700202
$ws.Cells.Item(25, 11).Value = 1.084298858483812
$ws.Cells.Item(25, 12).Value = 1.069006913159649
$ws.Cells.Item(25, 13).Value = 1.090436930055234
$ws.Cells.Item(25, 14).Value = 1.07898317263909
